# "correção nos dados e inicio da analise PNAD 2009"
#
# The row for the section header "grandes regiões e unidades da
# federação" (row 6) was a label-only row with no data beneath it.
# It is removed entirely: deleting the whole row shifts every row
# below it (the region/state data rows, through the last row for
# "goiás") up by one, so the region labels now line up directly with
# their data, and the shared-string entry for the removed header text
# is dropped along with it.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(6).Delete()
